$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 39 (1-based) / S/N "35", Date "16/6/2014" is the log entry being filled in.

# --- Name column (column 2): add the two description runs ---
$nameCell = $t.Cell(39, 2)
$nameRange = $nameCell.Range
$nameXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00CD51A0" w:rsidRPr="00CD51A0" w:rsidRDefault="00CD51A0" w:rsidP="00D30636"><w:pPr><w:pStyle w:val="TableStyle1"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cs="Helvetica"/><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cs="Helvetica"/><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>(Yanhao) Code in JAVA</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cs="Helvetica"/><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>(line creep spawn implementation, map mesh)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$nameRange.InsertXML($nameXml)

# --- Yanhao (hrs.) column (column 5): add the "8" hours run ---
$hrsCell = $t.Cell(39, 5)
$hrsRange = $hrsCell.Range
$hrsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00CD51A0" w:rsidRPr="00CD51A0" w:rsidRDefault="00CD51A0"><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Unicode MS" w:hAnsiTheme="minorHAnsi" w:cs="Helvetica"/><w:color w:val="000000"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Unicode MS" w:hAnsiTheme="minorHAnsi" w:cs="Helvetica"/><w:color w:val="000000"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>8</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$hrsRange.InsertXML($hrsXml)
